$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.716.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.571.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.00%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.59%  "
$ws.Range("E9").Value = "  -2.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.67%  "
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.348"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.032.66"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "62.595.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.81%  "
$ws.Range("E16").Value = "  -2.85%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.584.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "340.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("E20").Value = "  -2.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.57%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("E24").Value = "  -4.21%  "
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.164"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.00%  "
$ws.Range("B26").Value = "SuiNetwork"
$ws.Range("C26").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.30%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.56%  "
$ws.Range("E29").Value = "  -3.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "453.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0793"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.69%  "
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "176.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.58%  "
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.398"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.95%  "
$ws.Range("E37").Value = "  -2.31%  "
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("E40").Value = "  -3.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "39.99"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "156.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.631"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.03"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0531"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.75%  "
$ws.Range("E47").Value = "  -1.92%  "
$ws.Range("E48").Value = "  -3.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.76%  "
$ws.Range("E50").Value = "  +0.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.37%  "
